# Saldo.xlsx - reorder the first three data rows of the "Export" sheet.
#
# Before (rows 2-4):
#   004385806 ANILSON   108615.2
#   004482102 NATALIA    37567.8
#   005002457 ROSANGELA  34484.08
#
# After (rows 2-4):
#   004482102 NATALIA    37567.8
#   005002457 ROSANGELA  34484.08
#   004385806 ANILSON      4615.2   <- balance also updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read the current contents of rows 2-4 (A:C) BEFORE writing anything,
# using the explicit getter call - `.Value` as a bare property does not
# invoke on this host, `.Value()` does.
$a2 = $ws.Range("A2").Value()
$b2 = $ws.Range("B2").Value()

$a3 = $ws.Range("A3").Value()
$b3 = $ws.Range("B3").Value()
$c3 = $ws.Range("C3").Value()

$a4 = $ws.Range("A4").Value()
$b4 = $ws.Range("B4").Value()
$c4 = $ws.Range("C4").Value()

# Column A holds zero-padded account numbers ("004385806", ...). Force the
# destination cells to text first so the write doesn't re-infer them as
# numbers and strip the leading zeros.
$ws.Range("A2:A4").NumberFormat = "@"

# Row 2 <- old row 3 (NATALIA)
$ws.Range("A2").Value = $a3
$ws.Range("B2").Value = $b3
$ws.Range("C2").Value = $c3

# Row 3 <- old row 4 (ROSANGELA)
$ws.Range("A3").Value = $a4
$ws.Range("B3").Value = $b4
$ws.Range("C3").Value = $c4

# Row 4 <- old row 2 (ANILSON), with the updated balance
$ws.Range("A4").Value = $a2
$ws.Range("B4").Value = $b2
$ws.Range("C4").Value = 4615.2

# Drop the scratch "@" text format again so the cells end up with the same
# (default) styling they started with.
$ws.Range("A2:A4").ClearFormats()
